$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column CF: header date (copy CE1's style, then set new date value) ---
$ws.Range("CE1").Copy()
$ws.Range("CF1").PasteSpecial(-4122)
$ws.Range("CF1").Value = 45986

# --- CF3:CF138 = copy of CE3:CE138 (values only, no CE2/CE139 involved here) ---
$ws.Range("CE3:CE138").Copy()
$ws.Range("CF3:CF138").PasteSpecial(-4163)

# --- Row 139: CF139 gets its own (different) value ---
$ws.Range("CF139").Value = -0.2099036351493167

# --- Row 140: new CF140 value ---
$ws.Range("CF140").Value = 0

# --- Row 141: brand-new row, only A141 populated (copy A140's style) ---
$ws.Range("A140").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("A141").Value = 45976

$excel.CutCopyMode = 0
